$d = $word.ActiveDocument
$t = $d.Tables(1)

$t.Cell(1, 1).Range.Text = "31÷9=3, 4"
$t.Cell(1, 2).Range.Text = "82÷4=20, 2"
$t.Cell(1, 3).Range.Text = "19÷9=2, 1"
$t.Cell(1, 4).Range.Text = "56÷8=7, 0"
$t.Cell(1, 5).Range.Text = "93÷9=10, 3"
$t.Cell(5, 1).Range.Text = "64÷5=12, 4"
$t.Cell(5, 2).Range.Text = "75÷6=12, 3"
$t.Cell(5, 3).Range.Text = "40÷5=8, 0"
$t.Cell(5, 4).Range.Text = "23÷6=3, 5"
$t.Cell(5, 5).Range.Text = "49÷3=16, 1"
$t.Cell(9, 1).Range.Text = "35÷3=11, 2"
$t.Cell(9, 2).Range.Text = "42÷5=8, 2"
$t.Cell(9, 3).Range.Text = "28÷9=3, 1"
$t.Cell(9, 4).Range.Text = "62÷2=31, 0"
$t.Cell(9, 5).Range.Text = "90÷6=15, 0"
$t.Cell(13, 1).Range.Text = "92÷6=15, 2"
$t.Cell(13, 2).Range.Text = "42÷4=10, 2"
$t.Cell(13, 3).Range.Text = "28÷7=4, 0"
$t.Cell(13, 4).Range.Text = "81÷4=20, 1"
$t.Cell(13, 5).Range.Text = "44÷7=6, 2"
$t.Cell(17, 1).Range.Text = "93÷6=15, 3"
$t.Cell(17, 2).Range.Text = "38÷6=6, 2"
$t.Cell(17, 3).Range.Text = "85÷8=10, 5"
$t.Cell(17, 4).Range.Text = "22÷2=11, 0"
$t.Cell(17, 5).Range.Text = "25÷9=2, 7"
